$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target values for columns B..M, rows 1..15 (row 1 is the header row)
$data = @(
    @(1,2,3,4,5,6,7,8,9,10,11,12),
    @(0,0,0,0,0,0,73,0,0,0,15,15),
    @(43,0,0,0,0,0,18,12,6,26,0,24),
    @(111,0,0,0,0,183,0,0,0,0,0,3),
    @(0,0,0,47,1,0,0,75,0,0,33,27),
    @(0,0,0,0,3,0,16,0,0,20,0,0),
    @(10,0,0,36,0,0,0,28,0,0,0,12),
    @(0,0,0,0,0,0,0,14,0,0,0,3),
    @(2,0,0,4,14,0,0,23,0,0,20,0),
    @(25,0,15,25,63,57,45,0,0,95,0,60),
    @(48,0,120,660,0,120,120,270,0,0,570,0),
    @(10,0,0,12,174,0,0,0,12,72,0,60),
    @(0,0,0,41,0,0,45,100,0,35,35,55),
    @(1715,0,0,981,880,770,1760,0,440,330,770,550),
    @(21,0,5,30,0,25,85,0,20,15,45,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 1
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowValues[$j]
    }
}

# New columns D..M on the header row (row 1) need the same formatting as B1/C1
# (bold, centered/top-aligned, thin border). Copy the format from C1 and paste
# it (formats only) onto the newly populated header cells so the same style
# index gets reused instead of creating a brand-new style entry.
$ws.Range("C1").Copy()
$ws.Range("D1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
